# Estadisticos Segundo Parcial 26 Mayo
#
# Updates grades for "GESTIONA LOS PROCESOS DE CAPACITACION PARA EL
# DESARROLLO DEL TALENTO HUMANO" (1P -> column K, 3P -> column Y) on the
# "Calificaciones" sheet, and the corresponding attendance / average
# figures this feeds on "Asistencias" (columns K/R) and "Totales" (H4).

$wb = $excel.ActiveWorkbook

# --- Calificaciones: 1P (K) / 3P (Y) grades for "GESTIONA..." ---
$ws = $wb.Worksheets.Item("Calificaciones")
$ws.Range("K4").Value2  = 9
$ws.Range("K5").Value2  = 9
$ws.Range("Y5").Value2  = 9
$ws.Range("K6").Value2  = 9
$ws.Range("Y6").Value2  = 8
$ws.Range("K7").Value2  = 7
$ws.Range("K8").Value2  = 9
$ws.Range("Y8").Value2  = 9
$ws.Range("K9").Value2  = 8
$ws.Range("K10").Value2 = 9
$ws.Range("Y10").Value2 = 9
$ws.Range("K11").Value2 = 5
$ws.Range("K12").Value2 = 9
$ws.Range("K13").Value2 = 8
$ws.Range("K14").Value2 = 8
$ws.Range("K15").Value2 = 8
$ws.Range("K16").Value2 = 9
$ws.Range("K17").Value2 = 9
$ws.Range("K18").Value2 = 8
$ws.Range("Y18").Value2 = 7

# --- Asistencias: attendance % recalculated for the same subject/students ---
$ws = $wb.Worksheets.Item("Asistencias")
$ws.Range("K4").Value2  = 96.90000000000001
$ws.Range("R4").Value2  = 96.90000000000001
$ws.Range("K5").Value2  = 97.90000000000001
$ws.Range("R5").Value2  = 97.90000000000001
$ws.Range("K7").Value2  = 97.90000000000001
$ws.Range("R7").Value2  = 97.90000000000001
$ws.Range("K8").Value2  = 97.90000000000001
$ws.Range("R8").Value2  = 97.90000000000001
$ws.Range("K9").Value2  = 94.8
$ws.Range("R9").Value2  = 94.8
$ws.Range("K10").Value2 = 96.90000000000001
$ws.Range("R10").Value2 = 96.90000000000001
$ws.Range("K11").Value2 = 93.8
$ws.Range("R11").Value2 = 93.8
$ws.Range("K12").Value2 = 96.90000000000001
$ws.Range("R12").Value2 = 96.90000000000001
$ws.Range("K13").Value2 = 96.90000000000001
$ws.Range("R13").Value2 = 96.90000000000001
$ws.Range("K14").Value2 = 92.8
$ws.Range("R14").Value2 = 92.8
$ws.Range("K15").Value2 = 96.90000000000001
$ws.Range("R15").Value2 = 96.90000000000001
$ws.Range("K16").Value2 = 96.90000000000001
$ws.Range("R16").Value2 = 96.90000000000001
$ws.Range("K17").Value2 = 97.90000000000001
$ws.Range("R17").Value2 = 97.90000000000001
$ws.Range("K18").Value2 = 99
$ws.Range("R18").Value2 = 99

# --- Totales: group average for "GESTIONA..." ---
$ws = $wb.Worksheets.Item("Totales")
$ws.Range("H4").Value2 = 8.199999999999999
